$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K4").Value = 106.25
$ws.Range("H4").Value = 106.25
$ws.Range("M4").Value = 7.75
$ws.Range("I4").Value = 106.25
$ws.Range("H62").Value = 9748.375
$ws.Range("I62").Value = 13397.6
$ws.Range("M62").Value = -12773.6
$ws.Range("K62").Value = 13397.6
$ws.Range("M65").Value = -63868
$ws.Range("H65").Value = 9748.375
$ws.Range("K65").Value = 66988
$ws.Range("I65").Value = 13397.6
$ws.Range("H69").Value = 5499.5
$ws.Range("I69").Value = 5499.5
$ws.Range("K69").Value = 16498.5
$ws.Range("M69").Value = -15624.5
$ws.Range("H72").Value = 5499.5
$ws.Range("M72").Value = -45127.5
$ws.Range("K72").Value = 49495.5
$ws.Range("I72").Value = 5499.5
$ws.Range("M98").Value = 890.3125
$ws.Range("K98").Value = 607.6875
$ws.Range("H98").Value = 607.6875
$ws.Range("I98").Value = 607.6875
$ws.Range("I122").Value = 607.6875
$ws.Range("K122").Value = 1823.0625
$ws.Range("H122").Value = 607.6875
$ws.Range("M122").Value = 626.9375
$ws.Range("H137").Value = 4249.75
$ws.Range("M137").Value = -6448.5
$ws.Range("K137").Value = 8998.5
$ws.Range("I137").Value = 2999.5
$ws.Range("L138").Value = 16933.6362
$ws.Range("I138").Value = 2489.2222
$ws.Range("M138").Value = -2327.6666
$ws.Range("K138").Value = 7467.6666
$ws.Range("H138").Value = 4224.65
$ws.Range("J138").Value = 5644.5454
$ws.Range("N138").Value = -27213.6362
$ws.Range("N141").Value = -33592.49950000001
$ws.Range("L141").Value = 23232.4995
$ws.Range("I141").Value = 6022.0586
$ws.Range("K141").Value = 18066.1758
$ws.Range("M141").Value = -12886.1758
$ws.Range("J141").Value = 7744.1665
$ws.Range("H141").Value = 6471.304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("K32").Value = 4724.213
$ws.Range("H32").Value = 5144.1836
$ws.Range("M32").Value = -4437.213
$ws.Range("I32").Value = 4724.213
$ws.Range("K45").Value = 203298.6
$ws.Range("L45").Value = 8750
$ws.Range("J45").Value = 8750
$ws.Range("N45").Value = -9504
$ws.Range("H45").Value = 147713.28
$ws.Range("M45").Value = -202921.6
$ws.Range("I45").Value = 203298.6
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = ""
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("I122").Value = 2353.0833
$ws.Range("K122").Value = 7059.249899999999
$ws.Range("H122").Value = 2731.5
$ws.Range("M122").Value = -4609.249899999999
$ws.Range("N132").Value = -17195
$ws.Range("L132").Value = 12135
$ws.Range("J132").Value = 4045
$ws.Range("H132").Value = 1362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("J9").Value = 25000
$ws.Range("N9").Value = -25336
$ws.Range("N134").Value = -24313.5
$ws.Range("H134").Value = 3423.422
$ws.Range("L134").Value = 19243.5
$ws.Range("J134").Value = 6414.5
$ws.Range("L140").Value = 96000
$ws.Range("H140").Value = 96000
$ws.Range("N140").Value = -106360
$ws.Range("J140").Value = 96000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3077.7778
$ws.Range("J16").Value = 4438.8335
$ws.Range("L16").Value = 4438.8335
$ws.Range("N16").Value = -5012.8335
$ws.Range("K31").Value = 18376.777
$ws.Range("I31").Value = 18376.777
$ws.Range("M31").Value = -18081.777
$ws.Range("H31").Value = 10658.385
$ws.Range("H34").Value = 10658.385
$ws.Range("K34").Value = 18376.777
$ws.Range("I34").Value = 18376.777
$ws.Range("M34").Value = -18174.777
$ws.Range("K58").Value = 1049.75
$ws.Range("I58").Value = 1049.75
$ws.Range("H58").Value = 7729.5
$ws.Range("N58").Value = -10807.4
$ws.Range("L58").Value = 10401.4
$ws.Range("M58").Value = -846.75
$ws.Range("J58").Value = 10401.4
$ws.Range("H62").Value = 9999.666999999999
$ws.Range("I62").Value = 9999.666999999999
$ws.Range("M62").Value = -9375.666999999999
$ws.Range("K62").Value = 9999.666999999999
$ws.Range("M65").Value = -46878.335
$ws.Range("H65").Value = 9999.666999999999
$ws.Range("K65").Value = 49998.335
$ws.Range("I65").Value = 9999.666999999999
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("H86").Value = 9249.5
$ws.Range("M86").Value = ""
$ws.Range("M89").Value = ""
$ws.Range("I89").Value = 0
$ws.Range("H89").Value = 9249.5
$ws.Range("K89").Value = 0
$ws.Range("L97").Value = 30555
$ws.Range("J97").Value = 30555
$ws.Range("N97").Value = -32537
$ws.Range("H97").Value = 30555
$ws.Range("H113").Value = 3077.7778
$ws.Range("L113").Value = 4438.8335
$ws.Range("N113").Value = -8778.833500000001
$ws.Range("J113").Value = 4438.8335
$ws.Range("I132").Value = 3689.1
$ws.Range("M132").Value = -8537.299999999999
$ws.Range("K132").Value = 11067.3
$ws.Range("H132").Value = 4717.364
$ws.Range("K136").Value = 3149.25
$ws.Range("J136").Value = 10401.4
$ws.Range("L136").Value = 31204.2
$ws.Range("H136").Value = 7729.5
$ws.Range("N136").Value = -36304.2
$ws.Range("I136").Value = 1049.75
$ws.Range("M136").Value = -599.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L34").Value = 1260
$ws.Range("H34").Value = 420
$ws.Range("J34").Value = 420
$ws.Range("N34").Value = -1428
$ws.Range("N39").Value = -5041.125
$ws.Range("I39").Value = 426
$ws.Range("J39").Value = 1484.375
$ws.Range("M39").Value = -984
$ws.Range("H39").Value = 1272.7
$ws.Range("L39").Value = 4453.125
$ws.Range("K39").Value = 1278
$ws.Range("M92").Value = 573
$ws.Range("K92").Value = 675
$ws.Range("I92").Value = 225
$ws.Range("H92").Value = 225
$ws.Range("K97").Value = 7749.999899999999
$ws.Range("M97").Value = -7253.999899999999
$ws.Range("I97").Value = 2583.3333
$ws.Range("H97").Value = 2270
$ws.Range("M109").Value = -91586.00199999999
$ws.Range("H109").Value = 9116.786
$ws.Range("K109").Value = 92626.00199999999
$ws.Range("I109").Value = 30875.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K80").Value = 2498.75
$ws.Range("M80").Value = -1500.75
$ws.Range("I80").Value = 2498.75
$ws.Range("H80").Value = 2399
$ws.Range("M83").Value = -7501.75
$ws.Range("H83").Value = 2399
$ws.Range("K83").Value = 12493.75
$ws.Range("I83").Value = 2498.75
$ws.Range("I122").Value = 3098.25
$ws.Range("K122").Value = 9294.75
$ws.Range("H122").Value = 3109.5557
$ws.Range("M122").Value = -6844.75
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("H133").Value = 146000
$ws.Range("M133").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 1675
$ws.Range("H22").Value = 1974.8334
$ws.Range("N22").Value = -2265
$ws.Range("L22").Value = 1675
$ws.Range("H27").Value = 1974.8334
$ws.Range("J27").Value = 1675
$ws.Range("N27").Value = -1889
$ws.Range("L27").Value = 1675
$ws.Range("L40").Value = 4142.857
$ws.Range("J40").Value = 4142.857
$ws.Range("K40").Value = 2726.3076
$ws.Range("M40").Value = -2590.3076
$ws.Range("I40").Value = 2726.3076
$ws.Range("N40").Value = -4414.857
$ws.Range("H40").Value = 3222.1
$ws.Range("L94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H94").Value = 0
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("H98").Value = 50000
$ws.Range("N98").Value = -55990
$ws.Range("I100").Value = 4424.625
$ws.Range("K100").Value = 4424.625
$ws.Range("J100").Value = 6500
$ws.Range("L100").Value = 6500
$ws.Range("N100").Value = -7582
$ws.Range("M100").Value = -3883.625
$ws.Range("H100").Value = 4655.222
$ws.Range("L130").Value = 74992.5
$ws.Range("N130").Value = -85032.5
$ws.Range("J130").Value = 74992.5
$ws.Range("H130").Value = 74992.5
$ws.Range("I132").Value = 16055.643
$ws.Range("M132").Value = -45636.929
$ws.Range("N132").Value = -20060
$ws.Range("L132").Value = 15000
$ws.Range("K132").Value = 48166.929
$ws.Range("J132").Value = 5000
$ws.Range("H132").Value = 14673.6875
$ws.Range("K136").Value = 17169
$ws.Range("H136").Value = 5810.577
$ws.Range("I136").Value = 5723
$ws.Range("M136").Value = -14619

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N101").Value = -33592
$ws.Range("H101").Value = 27102
$ws.Range("L101").Value = 27102
$ws.Range("J101").Value = 27102
$ws.Range("K136").Value = 6927.558000000001
$ws.Range("J136").Value = 3024.1667
$ws.Range("L136").Value = 9072.500100000001
$ws.Range("H136").Value = 2396.7346
$ws.Range("N136").Value = -14172.5001
$ws.Range("I136").Value = 2309.186
$ws.Range("M136").Value = -4377.558000000001
